# Auto-generated edit script: updates computed profit columns (H-N)
# across multiple sheets per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 27027500
$ws.Range("I28").Value = 33333744
$ws.Range("J28").Value = 735.4286
$ws.Range("K28").Value = 33333744
$ws.Range("L28").Value = 735.4286
$ws.Range("M28").Value = -33333259
$ws.Range("N28").Value = -1705.4286
$ws.Range("H33").Value = 997.4666999999999
$ws.Range("I33").Value = 1262.6
$ws.Range("J33").Value = 467.2
$ws.Range("K33").Value = 1262.6
$ws.Range("L33").Value = 467.2
$ws.Range("M33").Value = -1033.6
$ws.Range("N33").Value = -925.2
$ws.Range("H74").Value = 8774
$ws.Range("I74").Value = 7698.6665
$ws.Range("K74").Value = 7698.6665
$ws.Range("M74").Value = -6762.6665
$ws.Range("H77").Value = 8774
$ws.Range("I77").Value = 7698.6665
$ws.Range("K77").Value = 38493.3325
$ws.Range("M77").Value = -33813.3325
$ws.Range("H98").Value = 2398.8
$ws.Range("I98").Value = 1887.2222
$ws.Range("K98").Value = 1887.2222
$ws.Range("M98").Value = -389.2221999999999
$ws.Range("H122").Value = 2398.8
$ws.Range("I122").Value = 1887.2222
$ws.Range("K122").Value = 5661.6666
$ws.Range("M122").Value = -3211.6666
$ws.Range("H129").Value = 2203.389
$ws.Range("J129").Value = 3263
$ws.Range("L129").Value = 9789
$ws.Range("N129").Value = -19789
$ws.Range("H137").Value = 4085.5
$ws.Range("I137").Value = 4277.575
$ws.Range("J137").Value = 2805
$ws.Range("K137").Value = 12832.725
$ws.Range("L137").Value = 8415
$ws.Range("M137").Value = -10282.725
$ws.Range("N137").Value = -13515

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 29367210
$ws.Range("I2").Value = 66073156
$ws.Range("J2").Value = 2452.7
$ws.Range("K2").Value = 66073156
$ws.Range("L2").Value = 2452.7
$ws.Range("M2").Value = -66073043
$ws.Range("N2").Value = -2678.7
$ws.Range("H32").Value = 14620.359
$ws.Range("I32").Value = 11579.857
$ws.Range("J32").Value = 34130.25
$ws.Range("K32").Value = 11579.857
$ws.Range("L32").Value = 34130.25
$ws.Range("M32").Value = -11292.857
$ws.Range("N32").Value = -34704.25
$ws.Range("H110").Value = 2953.2954
$ws.Range("I110").Value = 3106.475
$ws.Range("K110").Value = 3106.475
$ws.Range("M110").Value = -1061.475
$ws.Range("H116").Value = 29367210
$ws.Range("I116").Value = 66073156
$ws.Range("J116").Value = 2452.7
$ws.Range("K116").Value = 66073156
$ws.Range("L116").Value = 2452.7
$ws.Range("M116").Value = -66070862
$ws.Range("N116").Value = -7040.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 29367210
$ws.Range("I3").Value = 66073156
$ws.Range("J3").Value = 2452.7
$ws.Range("K3").Value = 66073156
$ws.Range("L3").Value = 2452.7
$ws.Range("M3").Value = -66073042
$ws.Range("N3").Value = -2680.7
$ws.Range("H22").Value = 168.85715
$ws.Range("J22").Value = 152
$ws.Range("L22").Value = 152
$ws.Range("N22").Value = -498
$ws.Range("H105").Value = 2160.5806
$ws.Range("I105").Value = 2085.4827
$ws.Range("K105").Value = 2085.4827
$ws.Range("M105").Value = -338.4827
$ws.Range("H107").Value = 2010.9131
$ws.Range("I107").Value = 1730.2
$ws.Range("J107").Value = 3882.3333
$ws.Range("K107").Value = 1730.2
$ws.Range("L107").Value = 3882.3333
$ws.Range("M107").Value = 189.8
$ws.Range("N107").Value = -7722.3333
$ws.Range("H134").Value = 3546.9788
$ws.Range("I134").Value = 3166.244
$ws.Range("J134").Value = 6148.6665
$ws.Range("K134").Value = 9498.732
$ws.Range("L134").Value = 18445.9995
$ws.Range("M134").Value = -6963.732
$ws.Range("N134").Value = -23515.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4007.7307
$ws.Range("I31").Value = 3022.8125
$ws.Range("J31").Value = 5583.6
$ws.Range("K31").Value = 3022.8125
$ws.Range("L31").Value = 5583.6
$ws.Range("M31").Value = -2727.8125
$ws.Range("N31").Value = -6173.6
$ws.Range("H34").Value = 4007.7307
$ws.Range("I34").Value = 3022.8125
$ws.Range("J34").Value = 5583.6
$ws.Range("K34").Value = 3022.8125
$ws.Range("L34").Value = 5583.6
$ws.Range("M34").Value = -2820.8125
$ws.Range("N34").Value = -5987.6
$ws.Range("H52").Value = 93316.664
$ws.Range("J52").Value = 93316.664
$ws.Range("L52").Value = 93316.664
$ws.Range("N52").Value = -93904.664
$ws.Range("H58").Value = 1605.1305
$ws.Range("I58").Value = 1086.6
$ws.Range("J58").Value = 2577.375
$ws.Range("K58").Value = 1086.6
$ws.Range("L58").Value = 2577.375
$ws.Range("M58").Value = -883.5999999999999
$ws.Range("N58").Value = -2983.375
$ws.Range("H99").Value = 9477.083000000001
$ws.Range("I99").Value = 6051.0435
$ws.Range("J99").Value = 15538.538
$ws.Range("K99").Value = 6051.0435
$ws.Range("L99").Value = 15538.538
$ws.Range("M99").Value = -4553.0435
$ws.Range("N99").Value = -18534.538
$ws.Range("H107").Value = 651.1818
$ws.Range("I107").Value = 616.3
$ws.Range("K107").Value = 616.3
$ws.Range("M107").Value = 1303.7
$ws.Range("H126").Value = 9477.083000000001
$ws.Range("I126").Value = 6051.0435
$ws.Range("J126").Value = 15538.538
$ws.Range("K126").Value = 18153.1305
$ws.Range("L126").Value = 46615.614
$ws.Range("M126").Value = -15683.1305
$ws.Range("N126").Value = -51555.614
$ws.Range("H136").Value = 1605.1305
$ws.Range("I136").Value = 1086.6
$ws.Range("J136").Value = 2577.375
$ws.Range("K136").Value = 3259.8
$ws.Range("L136").Value = 7732.125
$ws.Range("M136").Value = -709.7999999999997
$ws.Range("N136").Value = -12832.125
$ws.Range("H141").Value = 254327.67
$ws.Range("I141").Value = 66975
$ws.Range("J141").Value = 291798.2
$ws.Range("K141").Value = 66975
$ws.Range("L141").Value = 291798.2
$ws.Range("M141").Value = -61795
$ws.Range("N141").Value = -302158.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4971.25
$ws.Range("I116").Value = 4995
$ws.Range("J116").Value = 4900
$ws.Range("K116").Value = 14985
$ws.Range("L116").Value = 14700
$ws.Range("M116").Value = -11543
$ws.Range("N116").Value = -21584
$ws.Range("H122").Value = 990.2222
$ws.Range("I122").Value = 1145.125
$ws.Range("J122").Value = 866.3
$ws.Range("K122").Value = 10306.125
$ws.Range("L122").Value = 7796.7
$ws.Range("M122").Value = -7856.125
$ws.Range("N122").Value = -12696.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 506.08823
$ws.Range("I97").Value = 464.15384
$ws.Range("J97").Value = 642.375
$ws.Range("K97").Value = 464.15384
$ws.Range("L97").Value = 642.375
$ws.Range("M97").Value = 31.84616
$ws.Range("N97").Value = -1634.375
$ws.Range("H132").Value = 4651.1665
$ws.Range("I132").Value = 4444.351
$ws.Range("K132").Value = 13333.053
$ws.Range("M132").Value = -10803.053

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 963.46875
$ws.Range("J46").Value = 1018.087
$ws.Range("L46").Value = 1018.087
$ws.Range("N46").Value = -1394.087
$ws.Range("H82").Value = 1916.5807
$ws.Range("I82").Value = 1376.75
$ws.Range("J82").Value = 2898.0908
$ws.Range("K82").Value = 1376.75
$ws.Range("L82").Value = 2898.0908
$ws.Range("M82").Value = -1015.75
$ws.Range("N82").Value = -3620.0908
$ws.Range("H85").Value = 1916.5807
$ws.Range("I85").Value = 1376.75
$ws.Range("J85").Value = 2898.0908
$ws.Range("K85").Value = 1376.75
$ws.Range("L85").Value = 2898.0908
$ws.Range("M85").Value = -128.75
$ws.Range("N85").Value = -5394.0908
$ws.Range("H100").Value = 55558650
$ws.Range("I100").Value = 166669100
$ws.Range("K100").Value = 166669100
$ws.Range("M100").Value = -166668559

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6102.615
$ws.Range("I132").Value = 6547.5625
$ws.Range("J132").Value = 5390.7
$ws.Range("K132").Value = 19642.6875
$ws.Range("L132").Value = 16172.1
$ws.Range("M132").Value = -17112.6875
$ws.Range("N132").Value = -21232.1
$ws.Range("H137").Value = 149864.75
$ws.Range("J137").Value = 149864.75
$ws.Range("L137").Value = 149864.75
$ws.Range("N137").Value = -160064.75
